# Update "想去人数" (number of people interested) counts in column F
# for the "展览" (sheet1) and "全部类型" (sheet4) worksheets, matching the
# refreshed data snapshot pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")

$expoUpdates = @{
    3  = 14723
    4  = 18009
    5  = 18009
    16 = 63
    17 = 172
    19 = 1354
    20 = 151
    21 = 80
    22 = 71
    23 = 217
    24 = 7448
    25 = 983
    26 = 12
    28 = 1192
    30 = 5892
    31 = 76
    34 = 151
    35 = 240
    36 = 5169
}

foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")

$allUpdates = @{
    3  = 14723
    4  = 18009
    5  = 18009
    16 = 63
    17 = 172
    19 = 1354
    20 = 151
    21 = 80
    23 = 71
    24 = 217
    25 = 7448
    26 = 983
    27 = 12
    29 = 1192
    32 = 5892
    33 = 76
    36 = 151
    37 = 240
    38 = 5169
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}

$wb.Save()
